# Daily attendance processing - 2025-12-07 10:26:36
# Normalize the "Recorded By" (column G) entries on the "Session Analysis
# Results" sheet: the author/system identity that reflects the LAST
# recorder of a session should be listed FIRST, so for each affected row
# the first and last comma-separated names/emails in column G are swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column G whose "Recorded By" value needs the first and last
# comma-separated entries swapped.
$rows = @(2,3,5,6,8,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,31,32,34,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,57,58,60,62,63,64,65,66,67,69,70,71,72,73,74,76,78,80,81,82,83,84,85,86,90,92,93,94,96,99,101,106,107,108,109,110,111,112,116,118,119,120,122,125,127,132,133,134,135,136,137,138,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $old = $cell.Value2
    $parts = $old -split ", "
    if ($parts.Length -ge 2) {
        $first = $parts[0]
        $last = $parts[$parts.Length - 1]
        $parts[0] = $last
        $parts[$parts.Length - 1] = $first
        $new = $parts -join ", "
        $cell.Value = $new
    }
}
